$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.835941000000001
$ws.Range("H2").Value = 17.507823
$ws.Range("I2").Value = 0.03643643319117328
$ws.Range("J2").Value = 0.03643643319117327
$ws.Range("M2").Value = 10.34761366666667
$ws.Range("N2").Value = 31.042841
$ws.Range("O2").Value = 0.2299953477621856
$ws.Range("P2").Value = 0.2299953477621856
$ws.Range("Q2").Value = 60.38806284946034
$ws.Range("R2").Value = 543.492565645143
$ws.Range("S2").Value = 0.008380210123017538
$ws.Range("T2").Value = 0.008380210123017538

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.835941000000001
$ws.Range("H3").Value = 17.507823
$ws.Range("I3").Value = 0.03643643319117328
$ws.Range("J3").Value = 0.03643643319117327
$ws.Range("O3").Value = 0.6794731949692173
$ws.Range("P3").Value = 0.6794731949692174
$ws.Range("Q3").Value = 178.4039129554557
$ws.Range("R3").Value = 1605.635216599101
$ws.Range("S3").Value = 0.02475757967368894
$ws.Range("T3").Value = 0.02475757967368894

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.835941000000001
$ws.Range("H4").Value = 17.507823
$ws.Range("I4").Value = 0.03643643319117328
$ws.Range("J4").Value = 0.03643643319117327
$ws.Range("M4").Value = 4.073058666666666
$ws.Range("N4").Value = 12.219176
$ws.Range("O4").Value = 0.09053145726859702
$ws.Range("P4").Value = 0.09053145726859703
$ws.Range("Q4").Value = 23.77013006820534
$ws.Range("R4").Value = 213.931170613848
$ws.Range("S4").Value = 0.003298643394466793
$ws.Range("T4").Value = 0.003298643394466793

$ws.Range("G5").Value = 17.50798033333334
$ws.Range("H5").Value = 52.52394100000001
$ws.Range("I5").Value = 0.1093102818770573
$ws.Range("J5").Value = 0.1093102818770573
$ws.Range("M5").Value = 10.34761366666667
$ws.Range("N5").Value = 31.042841
$ws.Range("O5").Value = 0.2299953477621856
$ws.Range("P5").Value = 0.2299953477621856
$ws.Range("Q5").Value = 181.1658165729312
$ws.Range("R5").Value = 1630.492349156381
$ws.Range("S5").Value = 0.02514085629429632
$ws.Range("T5").Value = 0.02514085629429633

$ws.Range("G6").Value = 17.50798033333334
$ws.Range("H6").Value = 52.52394100000001
$ws.Range("I6").Value = 0.1093102818770573
$ws.Range("J6").Value = 0.1093102818770573
$ws.Range("O6").Value = 0.6794731949692173
$ws.Range("P6").Value = 0.6794731949692174
$ws.Range("Q6").Value = 535.2165485247075
$ws.Range("R6").Value = 4816.948936722368
$ws.Range("S6").Value = 0.07427340646998985
$ws.Range("T6").Value = 0.07427340646998985

$ws.Range("G7").Value = 17.50798033333334
$ws.Range("H7").Value = 52.52394100000001
$ws.Range("I7").Value = 0.1093102818770573
$ws.Range("J7").Value = 0.1093102818770573
$ws.Range("M7").Value = 4.073058666666666
$ws.Range("N7").Value = 12.219176
$ws.Range("O7").Value = 0.09053145726859702
$ws.Range("P7").Value = 0.09053145726859703
$ws.Range("Q7").Value = 71.3110310325129
$ws.Range("R7").Value = 641.799279292616
$ws.Range("S7").Value = 0.009896019112771107
$ws.Range("T7").Value = 0.009896019112771107

$ws.Range("G8").Value = 136.8238143333333
$ws.Range("H8").Value = 410.471443
$ws.Range("I8").Value = 0.8542532849317694
$ws.Range("J8").Value = 0.8542532849317694
$ws.Range("M8").Value = 10.34761366666667
$ws.Range("N8").Value = 31.042841
$ws.Range("O8").Value = 0.2299953477621856
$ws.Range("P8").Value = 0.2299953477621856
$ws.Range("Q8").Value = 1415.799971121063
$ws.Range("R8").Value = 12742.19974008956
$ws.Range("S8").Value = 0.1964742813448717
$ws.Range("T8").Value = 0.1964742813448717

$ws.Range("G9").Value = 136.8238143333333
$ws.Range("H9").Value = 410.471443
$ws.Range("I9").Value = 0.8542532849317694
$ws.Range("J9").Value = 0.8542532849317694
$ws.Range("O9").Value = 0.6794731949692173
$ws.Range("P9").Value = 0.6794731949692174
$ws.Range("Q9").Value = 4182.68516809156
$ws.Range("R9").Value = 37644.16651282404
$ws.Range("S9").Value = 0.5804422088255385
$ws.Range("T9").Value = 0.5804422088255387

$ws.Range("G10").Value = 136.8238143333333
$ws.Range("H10").Value = 410.471443
$ws.Range("I10").Value = 0.8542532849317694
$ws.Range("J10").Value = 0.8542532849317694
$ws.Range("M10").Value = 4.073058666666666
$ws.Range("N10").Value = 12.219176
$ws.Range("O10").Value = 0.09053145726859702
$ws.Range("P10").Value = 0.09053145726859703
$ws.Range("Q10").Value = 557.2914227767742
$ws.Range("R10").Value = 5015.622804990968
$ws.Range("S10").Value = 0.07733679476135912
$ws.Range("T10").Value = 0.07733679476135913
